$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.014.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.88%  "

$ws.Range("D3").Value = "'1.943.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.88%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'227.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.24%  "

$ws.Range("D6").Value = "'0.588"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.85%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'52.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.67%  "

$ws.Range("D9").Value = "'0.360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.56%  "

$ws.Range("D10").Value = "'56.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("E11").Value = "  -7.01%  "

$ws.Range("E12").Value = "  -4.57%  "

$ws.Range("D13").Value = "'2.232.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.78%  "

$ws.Range("D14").Value = "'13.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.94%  "

$ws.Range("D15").Value = "'0.734"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.53%  "

$ws.Range("D16").Value = "'19.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.62%  "

$ws.Range("D17").Value = "'1.963.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.34%  "

$ws.Range("D18").Value = "'4.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.62%  "

$ws.Range("D19").Value = "'35.948.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.73%  "

$ws.Range("D20").Value = "'66.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("E21").Value = "  -7.53%  "

$ws.Range("E22").Value = "  -5.68%  "

$ws.Range("D23").Value = "'219.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.83%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -12.77%  "

$ws.Range("D27").Value = "'159.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("D28").Value = "'8.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.85%  "

$ws.Range("D29").Value = "'18.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.50%  "

$ws.Range("E30").Value = "  -6.85%  "

$ws.Range("E31").Value = "  -11.20%  "

$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("D33").Value = "'4.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.55%  "

$ws.Range("D34").Value = "'0.0593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.69%  "

$ws.Range("D35").Value = "'4.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.98%  "

$ws.Range("E36").Value = "  -7.88%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -2.69%  "

$ws.Range("D39").Value = "'3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.14%  "

$ws.Range("D40").Value = "'2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").Value = "'4.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.80%  "

$ws.Range("D42").Value = "'1.382.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "

$ws.Range("D43").Value = "'0.0195"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.88%  "

$ws.Range("D44").Value = "'0.0847"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.92%  "

$ws.Range("E45").Value = "  -12.86%  "

$ws.Range("D46").Value = "'85.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.78%  "

$ws.Range("D47").Value = "'0.960"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.73%  "

$ws.Range("D48").Value = "'2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("D49").Value = "'14.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.78%  "

$ws.Range("D50").Value = "'6.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.28%  "

$ws.Range("D51").Value = "'2.126.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.87%  "

